# "Generate Report for Handoff"
#
# The localization-status report is refreshed: file b.md has now been
# handed off for localization (status moves from "Handed back: in sync
# with en-US" to "Ready for handoff"), new handback artifacts are
# recorded for b.md in both the zh-cn and de-de target-language sheets,
# and an out-of-date-handback error message is attached to b.md's row.

$wb = $excel.ActiveWorkbook

$statusOld       = "Handed back: in sync with en-US"
$statusNew       = "Ready for handoff"
$newTimestamp    = "2016-08-28 18:37:14"
$errorDetail     = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a24139c76d081efc7e139110d84aa25ccf8e15b0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce1294a6eb8482890ec0273f49fa62f758f57412/e2e/b.md."

# ---------------------------------------------------------------
# Overview sheet: b.md row (row 3) gets the new status + timestamp
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $newTimestamp

# ---------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusNew
# Leading apostrophe forces Excel to store this as literal text rather than
# auto-coercing the recognised word "False" into a real Boolean cell; reset
# the style afterwards so the quote-prefix marker doesn't linger.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-28 18:37:10"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1640625

# ---------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $newTimestamp
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1640625
